# Projects.xlsx update — adds 5 new project rows (ABBY_WEB, GitHubApi,
# PoolClub (MySnooker), RoyalBlood (RoyalBlood Studios), SignalRWeb (SignalR))
# into the alphabetically-ordered list on Sheet1, renumbers the "No" column
# for the shifted rows, and updates the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Helper-ish pattern: the sheet lists one project per row with a blank
# spacer row in between (row 3, 5, 7, ... are data rows; 4, 6, 8, ... are
# empty / not present at all in the sheet). To insert a new project we
# insert two fresh rows, copy number formatting from the row that is about
# to be pushed down (so the new row matches the surrounding style: centered
# "No" column, wrapped-text text columns), then blank out the spacer row
# completely (no cells at all, matching the rest of the sheet) and write
# the new values.
#
# We work from the bottom of the sheet upward so that earlier (lower-row)
# insert points are unaffected by later inserts done further down.
# ---------------------------------------------------------------------------

# 5) insert "SignalRWeb (SignalR)" before original row 27 ("Sindal")
$ws.Rows("27:28").Insert()
$ws.Range("A29:E29").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A27").Clear()
$ws.Range("D27:E27").Clear()
$ws.Range("A28:E28").Clear()
$ws.Range("B27").Value = "SignalRWeb (SignalR)"
$ws.Range("C27").Value = "Asp.net (Web Application) With C# "

# 4) insert "RoyalBlood (RoyalBlood Studios)" before original row 25 ("Ruber R1")
$ws.Rows("25:26").Insert()
$ws.Range("A27:E27").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A25").Clear()
$ws.Range("D25:E25").Clear()
$ws.Range("A26:E26").Clear()
$ws.Range("B25").Value = "RoyalBlood (RoyalBlood Studios)"
$ws.Range("C25").Value = "Asp.net (Web Application) With C# "

# 3) insert "PoolClub (MySnooker)" before original row 21 ("RouteManagerOnline")
$ws.Rows("21:22").Insert()
$ws.Range("A23:E23").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A21").Clear()
$ws.Range("D21:E21").Clear()
$ws.Range("A22:E22").Clear()
$ws.Range("B21").Value = "PoolClub (MySnooker)"
$ws.Range("C21").Value = "Asp.net (Web Application) With c# "

# 2) insert "GitHubApi" before original row 9 ("IntranetPortal")
$ws.Rows("9:10").Insert()
$ws.Range("A11:E11").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A9").Clear()
$ws.Range("D9:E9").Clear()
$ws.Range("A10:E10").Clear()
$ws.Range("B9").Value = "GitHubApi"
$ws.Range("C9").Value = "Asp.net (Web Application) With C# "

# 1) insert "ABBY_WEB" before original row 3 ("BlackArt") — this one keeps
# its "No" value (1) and also has an "Other Specifications" entry in E.
$ws.Rows("3:4").Insert()
$ws.Range("A5:E5").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)
$ws.Range("D3").Clear()
$ws.Range("A4:E4").Clear()
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ABBY_WEB"
$ws.Range("C3").Value = "Asp.net (Web Application) With C# "
$ws.Range("E3").Value = "ABBY Flexi Capture Api "

# ---------------------------------------------------------------------------
# Renumber the "No" column (A) for every surviving original row now that the
# five new rows have shifted things down.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 1
$ws.Range("A7").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A13").Value = 4
$ws.Range("A15").Value = 5
$ws.Range("A17").Value = 6
$ws.Range("A19").Value = 7
$ws.Range("A21").Value = 8
$ws.Range("A23").Value = 9
$ws.Range("A27").Value = 10
$ws.Range("A29").Value = 11
$ws.Range("A33").Value = 12
$ws.Range("A37").Value = 13
$ws.Range("A39").Value = 14
$ws.Range("A41").Value = 15
$ws.Range("A43").Value = 16
$ws.Range("A45").Value = 17
$ws.Range("A47").Value = 18
$ws.Range("A49").Value = 19
$ws.Range("A51").Value = 20

# ---------------------------------------------------------------------------
# Restore the saved view: scrolled so row 28 is at the top, with C35 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("C35").Select()

Write-Output "done"
